$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 3.845168516212527
$ws.Range("E2").Value = 2.271247526478293
$ws.Range("C3").Value = 4.156138964540079
$ws.Range("E3").Value = 2.54816413516763
$ws.Range("C4").Value = 6.813059503123609
$ws.Range("E4").Value = 3.304848950886519
$ws.Range("C5").Value = 3.15086384235741
$ws.Range("E5").Value = 2.75227989254978
$ws.Range("C6").Value = -0.8100302238676749
$ws.Range("E6").Value = 1.758104467524002
$ws.Range("C7").Value = 1.544811553068626
$ws.Range("E7").Value = 2.687259419948318
$ws.Range("C8").Value = 2.745283889750105
$ws.Range("E8").Value = 2.766469441386121
$ws.Range("C9").Value = 0.6185769761447801
$ws.Range("E9").Value = 2.370470941984903
$ws.Range("C10").Value = 2.176690369274925
$ws.Range("E10").Value = 2.578250867688547
$ws.Range("C11").Value = 2.613312245287847
$ws.Range("E11").Value = 2.754147109475591
$ws.Range("C12").Value = 1.239638987083946
$ws.Range("E12").Value = 2.493426994223724
$ws.Range("C13").Value = 2.513051972228775
$ws.Range("E13").Value = 2.489529953686964
$ws.Range("C14").Value = 1.921237443390145
$ws.Range("E14").Value = 2.12762431017075
$ws.Range("C15").Value = 0.1585342120028033
$ws.Range("E15").Value = 1.634581066715657
$ws.Range("C16").Value = 0.4972865286980932
$ws.Range("E16").Value = 1.245013078987922
$ws.Range("C17").Value = 0.802283903038159
$ws.Range("E17").Value = 1.117253368207849
$ws.Range("C18").Value = 1.124230726200648
$ws.Range("E18").Value = 1.748979051430455
$ws.Range("C19").Value = 1.380216232465781
$ws.Range("E19").Value = 1.427113753179343
$ws.Range("C20").Value = 3.178515672370708
$ws.Range("E20").Value = 1.673213377067717
$ws.Range("C21").Value = 2.572874367914246
$ws.Range("E21").Value = 2.26445559922297
$ws.Range("C22").Value = -7.574105122480734
$ws.Range("E22").Value = -13.03389973386627
$ws.Range("C23").Value = 0.3088110144677358
$ws.Range("E23").Value = 0.8151717706756934
$ws.Range("C24").Value = 3.845022163334466
$ws.Range("E24").Value = 2.160819775353762
$ws.Range("C25").Value = 1.393976211508785
$ws.Range("E25").Value = 1.299377455687734
$ws.Range("C26").Value = -0.1149334425204152
$ws.Range("E26").Value = 0.7902221915852214
$ws.Range("C27").Value = 1.996585218152536
$ws.Range("E27").Value = 1.687142382558493
$ws.Range("C28").Value = 1.344646667329186
$ws.Range("E28").Value = 1.307762937264378
$ws.Range("C29").Value = 1.634452213710591
$ws.Range("E29").Value = 1.652249888306234
$ws.Range("C30").Value = 1.77159272734837
$ws.Range("E30").Value = 1.740794503770515
$ws.Range("C31").Value = 2.197572811344384
$ws.Range("E31").Value = 1.445892187720954
$ws.Range("C32").Value = 0.7727659837340584
$ws.Range("E32").Value = 1.351596145190048
$ws.Range("C33").Value = -1.9842724512803
$ws.Range("E33").Value = -1.032433245924502
$ws.Range("C34").Value = 0.3025215973878836
$ws.Range("E34").Value = 0.9864098545071087
$ws.Range("C35").Value = 2.147433528392373
$ws.Range("E35").Value = 1.587614009576455
$ws.Range("C36").Value = 0.2531386145756764
$ws.Range("E36").Value = 1.171952640618934
$ws.Range("C37").Value = 0.1066070725092017
$ws.Range("E37").Value = 1.106834879074126
$ws.Range("C38").Value = 0.3236660530042679
$ws.Range("E38").Value = 1.117873748384368
